# Auto-generated edit script applying the diff to Ifrit_Profits (market price refresh)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2527.8909
$ws.Range("J17").Value = 2527.8909
$ws.Range("L17").Value = 7583.672699999999
$ws.Range("N17").Value = -7919.672699999999
$ws.Range("H76").Value = 4005.0322
$ws.Range("I76").Value = 4267.6523
$ws.Range("J76").Value = 3250
$ws.Range("K76").Value = 4267.6523
$ws.Range("L76").Value = 3250
$ws.Range("M76").Value = -3952.6523
$ws.Range("N76").Value = -3880
$ws.Range("H79").Value = 4005.0322
$ws.Range("I79").Value = 4267.6523
$ws.Range("J79").Value = 3250
$ws.Range("K79").Value = 4267.6523
$ws.Range("L79").Value = 3250
$ws.Range("M79").Value = -3175.6523
$ws.Range("N79").Value = -5434
$ws.Range("H116").Value = 2042.8572
$ws.Range("I116").Value = 1883.3334
$ws.Range("K116").Value = 1883.3334
$ws.Range("M116").Value = 1558.6666
$ws.Range("H141").Value = 2276.2666
$ws.Range("I141").Value = 1724.5714
$ws.Range("K141").Value = 5173.7142
$ws.Range("M141").Value = 6.285799999999654

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9652.743
$ws.Range("I32").Value = 5182.603
$ws.Range("K32").Value = 5182.603
$ws.Range("M32").Value = -4895.603
$ws.Range("H61").Value = 4275084.5
$ws.Range("I61").Value = 5051957.5
$ws.Range("J61").Value = 2282
$ws.Range("K61").Value = 5051957.5
$ws.Range("L61").Value = 2282
$ws.Range("M61").Value = -5051745.5
$ws.Range("N61").Value = -2706
$ws.Range("H63").Value = 3867.3333
$ws.Range("I63").Value = 1950
$ws.Range("J63").Value = 5401.2
$ws.Range("K63").Value = 1950
$ws.Range("L63").Value = 5401.2
$ws.Range("M63").Value = -1264
$ws.Range("N63").Value = -6773.2
$ws.Range("H66").Value = 3867.3333
$ws.Range("I66").Value = 1950
$ws.Range("J66").Value = 5401.2
$ws.Range("K66").Value = 9750
$ws.Range("L66").Value = 27006
$ws.Range("M66").Value = -6318
$ws.Range("N66").Value = -33870
$ws.Range("H74").Value = 13047368
$ws.Range("I74").Value = 17647740
$ws.Range("K74").Value = 17647740
$ws.Range("M74").Value = -17646866
$ws.Range("H77").Value = 13047368
$ws.Range("I77").Value = 17647740
$ws.Range("K77").Value = 88238700
$ws.Range("M77").Value = -88234332
$ws.Range("H124").Value = 10166.5
$ws.Range("J124").Value = 10166.5
$ws.Range("L124").Value = 10166.5
$ws.Range("N124").Value = -19986.5
$ws.Range("H125").Value = 37000
$ws.Range("J125").Value = 37000
$ws.Range("L125").Value = 37000
$ws.Range("N125").Value = -46840
$ws.Range("H132").Value = 1078120.4
$ws.Range("I132").Value = 1570366.1
$ws.Range("J132").Value = 93628.91
$ws.Range("K132").Value = 4711098.300000001
$ws.Range("L132").Value = 280886.73
$ws.Range("M132").Value = -4708568.300000001
$ws.Range("N132").Value = -285946.73
$ws.Range("H136").Value = 4275084.5
$ws.Range("I136").Value = 5051957.5
$ws.Range("J136").Value = 2282
$ws.Range("K136").Value = 15155872.5
$ws.Range("L136").Value = 6846
$ws.Range("M136").Value = -15153322.5
$ws.Range("N136").Value = -11946

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1677.5
$ws.Range("I105").Value = 1568.125
$ws.Range("K105").Value = 1568.125
$ws.Range("M105").Value = 178.875
$ws.Range("H107").Value = 404181.3
$ws.Range("I107").Value = 564979.2
$ws.Range("J107").Value = 2186.6
$ws.Range("K107").Value = 564979.2
$ws.Range("L107").Value = 2186.6
$ws.Range("M107").Value = -563059.2
$ws.Range("N107").Value = -6026.6
$ws.Range("H134").Value = 22334500
$ws.Range("I134").Value = 22334500
$ws.Range("K134").Value = 67003500
$ws.Range("M134").Value = -67000965

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 606
$ws.Range("I17").Value = 545
$ws.Range("J17").Value = 850
$ws.Range("K17").Value = 1635
$ws.Range("L17").Value = 2550
$ws.Range("M17").Value = -1466
$ws.Range("N17").Value = -2888
$ws.Range("H25").Value = 2347.1428
$ws.Range("I25").Value = 1357.5
$ws.Range("J25").Value = 3666.6667
$ws.Range("K25").Value = 4072.5
$ws.Range("L25").Value = 11000.0001
$ws.Range("M25").Value = -3903.5
$ws.Range("N25").Value = -11338.0001
$ws.Range("H30").Value = 2347.1428
$ws.Range("I30").Value = 1357.5
$ws.Range("J30").Value = 3666.6667
$ws.Range("K30").Value = 4072.5
$ws.Range("L30").Value = 11000.0001
$ws.Range("M30").Value = -3970.5
$ws.Range("N30").Value = -11204.0001
$ws.Range("H34").Value = 2503.9048
$ws.Range("I34").Value = 149.75
$ws.Range("K34").Value = 449.25
$ws.Range("M34").Value = -365.25
$ws.Range("H39").Value = 11522.667
$ws.Range("J39").Value = 11522.667
$ws.Range("L39").Value = 34568.001
$ws.Range("N39").Value = -35156.001
$ws.Range("H55").Value = 1410.5264
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1410.5264
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 4231.5792
$ws.Range("N55").Value = -4585.5792
$ws.Range("H61").Value = 274.6
$ws.Range("I61").Value = 72.75
$ws.Range("J61").Value = 409.16666
$ws.Range("K61").Value = 218.25
$ws.Range("L61").Value = 1227.49998
$ws.Range("M61").Value = -3.25
$ws.Range("N61").Value = -1657.49998
$ws.Range("H70").Value = 91388.5
$ws.Range("I70").Value = 91388.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 274165.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -273850.5
$ws.Range("H73").Value = 91388.5
$ws.Range("I73").Value = 91388.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 274165.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -273073.5
$ws.Range("H113").Value = 500.84616
$ws.Range("I113").Value = 531.1429000000001
$ws.Range("J113").Value = 483.88
$ws.Range("K113").Value = 1593.4287
$ws.Range("L113").Value = 1451.64
$ws.Range("M113").Value = 576.5712999999998
$ws.Range("N113").Value = -5791.639999999999
$ws.Range("H122").Value = 14881803
$ws.Range("I122").Value = 19608402
$ws.Range("J122").Value = 3402918.5
$ws.Range("K122").Value = 176475618
$ws.Range("L122").Value = 30626266.5
$ws.Range("M122").Value = -176473168
$ws.Range("N122").Value = -30631166.5
$ws.Range("H131").Value = 736.33
$ws.Range("I131").Value = 527.3125
$ws.Range("J131").Value = 776.1429000000001
$ws.Range("K131").Value = 1581.9375
$ws.Range("L131").Value = 2328.4287
$ws.Range("M131").Value = 3458.0625
$ws.Range("N131").Value = -12408.4287
$ws.Range("M55").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5130.769
$ws.Range("I70").Value = 5360
$ws.Range("J70").Value = 4987.5
$ws.Range("K70").Value = 5360
$ws.Range("L70").Value = 4987.5
$ws.Range("M70").Value = -5090
$ws.Range("N70").Value = -5527.5
$ws.Range("H73").Value = 5130.769
$ws.Range("I73").Value = 5360
$ws.Range("J73").Value = 4987.5
$ws.Range("K73").Value = 5360
$ws.Range("L73").Value = 4987.5
$ws.Range("M73").Value = -4424
$ws.Range("N73").Value = -6859.5
$ws.Range("H102").Value = 2198.054
$ws.Range("I102").Value = 1967.2
$ws.Range("J102").Value = 2679
$ws.Range("K102").Value = 1967.2
$ws.Range("L102").Value = 2679
$ws.Range("M102").Value = -345.2
$ws.Range("N102").Value = -5923

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7950.55
$ws.Range("I132").Value = 9729.5
$ws.Range("J132").Value = 3799.6667
$ws.Range("K132").Value = 29188.5
$ws.Range("L132").Value = 11399.0001
$ws.Range("M132").Value = -26658.5
$ws.Range("N132").Value = -16459.0001
$ws.Range("H136").Value = 2490.8
$ws.Range("J136").Value = 5464.1665
$ws.Range("L136").Value = 16392.4995
$ws.Range("N136").Value = -21492.4995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 929.1142599999999
$ws.Range("I122").Value = 906.13043
$ws.Range("J122").Value = 973.1667
$ws.Range("K122").Value = 2718.39129
$ws.Range("L122").Value = 2919.5001
$ws.Range("M122").Value = -268.39129
$ws.Range("N122").Value = -7819.5001
